{"js": "// Replace the answer text in each populated cell of the division-problems\n// table with the new value, cell-by-cell (row, col) so that duplicate\n// strings (e.g. \"59\u00f77=8, 3\" appears twice) are handled positionally and\n// each run's formatting (rFonts/sz) is left untouched \u2014 only the <w:t>\n// text content changes, matching the source diff exactly.\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"28\u00f73=9, 1\", newText: \"66\u00f77=9, 3\" },\n  { row: 0, col: 1, oldText: \"99\u00f76=16, 3\", newText: \"29\u00f73=9, 2\" },\n  { row: 0, col: 2, oldText: \"42\u00f74=10, 2\", newText: \"14\u00f78=1, 6\" },\n  { row: 0, col: 3, oldText: \"33\u00f75=6, 3\", newText: \"52\u00f72=26, 0\" },\n  { row: 0, col: 4, oldText: \"96\u00f79=10, 6\", newText: \"67\u00f79=7, 4\" },\n\n  { row: 4, col: 0, oldText: \"49\u00f78=6, 1\", newText: \"10\u00f79=1, 1\" },\n  { row: 4, col: 1, oldText: \"22\u00f73=7, 1\", newText: \"38\u00f75=7, 3\" },\n  { row: 4, col: 2, oldText: \"89\u00f75=17, 4\", newText: \"10\u00f72=5, 0\" },\n  { row: 4, col: 3, oldText: \"14\u00f72=7, 0\", newText: \"51\u00f77=7, 2\" },\n  { row: 4, col: 4, oldText: \"91\u00f77=13, 0\", newText: \"48\u00f73=16, 0\" },\n\n  { row: 8, col: 0, oldText: \"13\u00f79=1, 4\", newText: \"21\u00f78=2, 5\" },\n  { row: 8, col: 1, oldText: \"55\u00f78=6, 7\", newText: \"63\u00f79=7, 0\" },\n  { row: 8, col: 2, oldText: \"17\u00f79=1, 8\", newText: \"75\u00f76=12, 3\" },\n  { row: 8, col: 3, oldText: \"27\u00f77=3, 6\", newText: \"62\u00f74=15, 2\" },\n  { row: 8, col: 4, oldText: \"64\u00f74=16, 0\", newText: \"42\u00f73=14, 0\" },\n\n  { row: 12, col: 0, oldText: \"59\u00f77=8, 3\", newText: \"26\u00f73=8, 2\" },\n  { row: 12, col: 1, oldText: \"93\u00f76=15, 3\", newText: \"71\u00f76=11, 5\" },\n  { row: 12, col: 2, oldText: \"42\u00f79=4, 6\", newText: \"70\u00f74=17, 2\" },\n  { row: 12, col: 3, oldText: \"77\u00f75=15, 2\", newText: \"49\u00f74=12, 1\" },\n  { row: 12, col: 4, oldText: \"63\u00f75=12, 3\", newText: \"44\u00f75=8, 4\" },\n\n  { row: 16, col: 0, oldText: \"48\u00f76=8, 0\", newText: \"19\u00f77=2, 5\" },\n  { row: 16, col: 1, oldText: \"93\u00f75=18, 3\", newText: \"81\u00f76=13, 3\" },\n  { row: 16, col: 2, oldText: \"26\u00f79=2, 8\", newText: \"29\u00f78=3, 5\" },\n  { row: 16, col: 3, oldText: \"59\u00f77=8, 3\", newText: \"35\u00f75=7, 0\" },\n  { row: 16, col: 4, oldText: \"67\u00f77=9, 4\", newText: \"63\u00f77=9, 0\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  cell.value = r.newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the answer text in each populated cell of the division-problems\n# table with the new value, cell-by-cell (Row, Col) so that duplicate\n# strings (e.g. \"59\u00f77=8, 3\" appears twice) are handled positionally and\n# each run's formatting (rFonts/sz) is left untouched - only the cell's\n# text content changes, matching the source diff exactly.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; NewText = \"66\u00f77=9, 3\" },\n    @{ Row = 1; Col = 2; NewText = \"29\u00f73=9, 2\" },\n    @{ Row = 1; Col = 3; NewText = \"14\u00f78=1, 6\" },\n    @{ Row = 1; Col = 4; NewText = \"52\u00f72=26, 0\" },\n    @{ Row = 1; Col = 5; NewText = \"67\u00f79=7, 4\" },\n\n    @{ Row = 5; Col = 1; NewText = \"10\u00f79=1, 1\" },\n    @{ Row = 5; Col = 2; NewText = \"38\u00f75=7, 3\" },\n    @{ Row = 5; Col = 3; NewText = \"10\u00f72=5, 0\" },\n    @{ Row = 5; Col = 4; NewText = \"51\u00f77=7, 2\" },\n    @{ Row = 5; Col = 5; NewText = \"48\u00f73=16, 0\" },\n\n    @{ Row = 9; Col = 1; NewText = \"21\u00f78=2, 5\" },\n    @{ Row = 9; Col = 2; NewText = \"63\u00f79=7, 0\" },\n    @{ Row = 9; Col = 3; NewText = \"75\u00f76=12, 3\" },\n    @{ Row = 9; Col = 4; NewText = \"62\u00f74=15, 2\" },\n    @{ Row = 9; Col = 5; NewText = \"42\u00f73=14, 0\" },\n\n    @{ Row = 13; Col = 1; NewText = \"26\u00f73=8, 2\" },\n    @{ Row = 13; Col = 2; NewText = \"71\u00f76=11, 5\" },\n    @{ Row = 13; Col = 3; NewText = \"70\u00f74=17, 2\" },\n    @{ Row = 13; Col = 4; NewText = \"49\u00f74=12, 1\" },\n    @{ Row = 13; Col = 5; NewText = \"44\u00f75=8, 4\" },\n\n    @{ Row = 17; Col = 1; NewText = \"19\u00f77=2, 5\" },\n    @{ Row = 17; Col = 2; NewText = \"81\u00f76=13, 3\" },\n    @{ Row = 17; Col = 3; NewText = \"29\u00f78=3, 5\" },\n    @{ Row = 17; Col = 4; NewText = \"35\u00f75=7, 0\" },\n    @{ Row = 17; Col = 5; NewText = \"63\u00f77=9, 0\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $table.Cell($r.Row, $r.Col)\n    $cell.Range.Text = $r.NewText\n}\n"}
